# Add 2022-05-20 data update to the Fonds de solidarite regional NAF dataset.
# Updates columns C (nombre_aides), D (nombre_entreprises) and E (montant_total)
# for the listed rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 9;   C = 69571;  D = 10018; E = 191359921 },
    @{ Row = 17;  C = 134741; D = 19602; E = 296792831 },
    @{ Row = 122; C = 9694;   D = 1445;  E = 31947576 },
    @{ Row = 164; C = 50573;  D = 11057; E = 168517372 },
    @{ Row = 168; C = 284958; D = 58122; E = 1209031755 },
    @{ Row = 169; C = 562583; D = 60954; E = 1284389559 },
    @{ Row = 170; C = 367327; D = 38110; E = 2844804389 },
    @{ Row = 171; C = 115128; D = 20264; E = 445409607 },
    @{ Row = 173; C = 54387;  D = 11601; E = 151858805 },
    @{ Row = 174; C = 357191; D = 69788; E = 1016884274 },
    @{ Row = 175; C = 125524; D = 18100; E = 812307519 },
    @{ Row = 179; C = 235683; D = 29335; E = 812488760 },
    @{ Row = 204; C = 4757;   D = 729;   E = 11756409 },
    @{ Row = 205; C = 11125;  D = 1314;  E = 44114295 },
    @{ Row = 209; C = 5364;   D = 1136;  E = 12211202 },
    @{ Row = 247; C = 29424;  D = 3751;  E = 99451778 },
    @{ Row = 264; C = 47471;  D = 7784;  E = 81946515 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
}

$wb.Save()
